$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns (match the bold/bordered header style
# used by the existing header row, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I column is 1 for every row, except rows 24 and 25 which mirror H (IP)
# J column mirrors H (IP) for every row, except rows 24 and 25 which differ

$iValues = @{
    2 = 1; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1; 20 = 1;
    21 = 1; 22 = 1; 23 = 1; 24 = 4; 25 = 3; 26 = 1; 27 = 1; 28 = 1; 29 = 1; 30 = 1;
    31 = 1; 32 = 1; 33 = 1; 34 = 1; 35 = 1; 36 = 1
}

$jValues = @{
    2 = 6; 3 = 5; 4 = 6; 5 = 6; 6 = 7; 7 = 6; 8 = 6; 9 = 4; 10 = 2;
    11 = 6; 12 = 6; 13 = 6; 14 = 6; 15 = 6; 16 = 7; 17 = 6; 18 = 5; 19 = 5; 20 = 6;
    21 = 6; 22 = 6; 23 = 5; 24 = 7; 25 = 5; 26 = 5; 27 = 5; 28 = 6; 29 = 5; 30 = 6;
    31 = 4; 32 = 4; 33 = 4; 34 = 4; 35 = 3; 36 = 2
}

foreach ($r in 2..36) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
